# Update "想去人数" (F column) figures on both the "展览" sheet and the
# "全部类型" sheet (which mirrors the same data) as of the latest scrape.

$wb = $excel.ActiveWorkbook

$updates = @{
    2 = 1344
    3 = 1902
    4 = 185
    6 = 6321
    7 = 187
    8 = 112
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}
